$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: the "Typical Limit Switch / Potentiometer" block now also
# documents Encoder A/B signaling, so update its title.
# ---------------------------------------------------------------------
$ws.Range("S4").Value = "Typical Limit Switch / Potentiometer / Encoder"

# ---------------------------------------------------------------------
# Pin 1 row (row 7): "A Chan" -> "Up" for the limit-switch style wiring,
# and the per-pin Color column for the encoder tables is cleared since
# that information no longer applies once the channel is generalized.
# ---------------------------------------------------------------------
$ws.Range("H7").Value = "Up"
$ws.Range("I7").ClearContents()
$ws.Range("L7").Value = "Up"
$ws.Range("M7").ClearContents()
$ws.Range("T7").Value = "Up"

# ---------------------------------------------------------------------
# Pin 2 row (row 8): "B Chan" -> "Down"
# ---------------------------------------------------------------------
$ws.Range("H8").Value = "Down"
$ws.Range("I8").ClearContents()
$ws.Range("L8").Value = "Down"
$ws.Range("M8").ClearContents()
$ws.Range("T8").Value = "Down"

# ---------------------------------------------------------------------
# Pin 3 row (row 9): encoders gain a Position (AIO) signal; the limit
# switch table's existing "position" label is re-cased to "Position".
# ---------------------------------------------------------------------
$ws.Range("G9").Value = "AIO-n (sig)"
$ws.Range("H9").Value = "Position"
$ws.Range("K9").Value = "AIO-n (sig)"
$ws.Range("L9").Value = "Position"
$ws.Range("T9").Value = "Position"

# ---------------------------------------------------------------------
# Pin 4 row (row 10): new "A" channel entries for both encoders and the
# limit switch / potentiometer / encoder table.
# ---------------------------------------------------------------------
$ws.Range("G10").Value = "DIO-n (sig)"
$ws.Range("H10").Value = "A"
$ws.Range("I10").Value = "Blue"
$ws.Range("K10").Value = "DIO-n (sig)"
$ws.Range("L10").Value = "A"
$ws.Range("M10").Value = "Blue"
$ws.Range("S10").Value = "DIO-n (sig)"
$ws.Range("T10").Value = "A"
$ws.Range("L10:M10").HorizontalAlignment = -4108
$ws.Range("S10:T10").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Pin 5 row (row 11): new "B" channel entries for both encoders and the
# limit switch / potentiometer / encoder table.
# ---------------------------------------------------------------------
$ws.Range("G11").Value = "DIO-n (sig)"
$ws.Range("H11").Value = "B"
$ws.Range("I11").Value = "Yellow"
$ws.Range("K11").Value = "DIO-n (sig)"
$ws.Range("L11").Value = "B"
$ws.Range("M11").Value = "Yellow"
$ws.Range("S11").Value = "DIO-n (sig)"
$ws.Range("T11").Value = "B"

# ---------------------------------------------------------------------
# Restore the selection as left by the editor.
# ---------------------------------------------------------------------
$ws.Range("M23").Select()
